$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.591.53'
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").Value = '2.991.13'
$ws.Range("E3").Value = '  +3.15%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '381.04'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.02'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.546'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.01%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.597'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.45'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.52%  '
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0849'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.69%  '
$ws.Range("D13").Value = '3.456.93'
$ws.Range("E13").Value = '  +3.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.48'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.60'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.93%  '
$ws.Range("D16").Value = '2.986.36'
$ws.Range("E16").Value = '  +3.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.978'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +6.82%  '
$ws.Range("D18").Value = '51.543.82'
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.34'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.46'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.04'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.48%  '
$ws.Range("D22").Value = '0.0₃0966'
$ws.Range("E22").Value = '  +3.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.43'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '263.21'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.90'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +9.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.34'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +20.31%  '
$ws.Range("E27").Value = '  +27.84%  '
$ws.Range("E28").Value = '  +16.89%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.172'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '26.08'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.64%  '
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.93'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.24'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +4.23%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '51.09'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("B35").Value = 'Toncoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.08'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0455'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +8.81%  '
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.06'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.23'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.59'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.86'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.117'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +4.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '125.76'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +6.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.02'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.284'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +21.78%  '
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.38'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.63%  '
$ws.Range("D48").Value = '2.041.81'
$ws.Range("E48").Value = '  +1.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.28'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +5.17%  '
$ws.Range("E50").Value = '  +8.69%  '
$ws.Range("E51").Value = '  +4.67%  '
